# aggiornamento fino a 21 marzo
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) of the last existing data row (229) down into the
# four new rows so the new date cells (column A) keep the same date style.
$ws.Range("A229:D229").Copy() | Out-Null
$ws.Range("A230:D233").PasteSpecial(-4122) | Out-Null

$newRows = @(
    @(44304, 3, 7, 106.6098081023454),
    @(44305, 0, 7, 106.6098081023454),
    @(44306, 0, 7, 106.6098081023454),
    @(44307, 1, 8, 121.8397806883948)
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = 230 + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newRows[$i][2]
    $ws.Cells.Item($r, 4).Value = $newRows[$i][3]
}
